# Added change user pass.
#
# The "8. Kao korisnik potrebno je da mogu da promenim svoj PASS" user
# story (row 58) has its task breakdown filled in with the actual time
# spent ("Realno utroseno vreme", column C) for rows 60-62, and the
# "10. Sistem treba da racuna proviziju..." story's row 74 also gets its
# actual time recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: scroll the view so row 52 is at the top-left, matching the
# author's on-screen position when they made this edit. Not all hosts
# persist window/scroll state to the saved file, so this is wrapped
# defensively and never blocks the data edits below.
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 52
    $win.ScrollColumn = 1
} catch {
}

# Fill in the "Realno utroseno vreme" (actual time spent) values that were
# previously left blank.
$ws.Range("C60").Value = "5min"
$ws.Range("C61").Value = "10min"
$ws.Range("C62").Value = "5min"
$ws.Range("C74").Value = "5min"
